$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 65, shifting existing rows 65-108 down to 66-109.
$ws.Rows.Item(65).Insert()

# Populate the newly inserted row 65 with its data.
$ws.Range("A65").Value() = 1
$ws.Range("B65").Value() = "Agrícola del Norte S.A. de Arica"
$ws.Range("C65").Value() = "Arica y Parinacota"
$ws.Range("D65").Value() = 44762
$ws.Range("E65").Value() = 15
$ws.Range("F65").Value() = "Fruta"
$ws.Range("G65").Value() = 100102
$ws.Range("H65").Value() = "Cítricos"
$ws.Range("I65").Value() = 100102004
$ws.Range("J65").Value() = "Mandarina"
$ws.Range("K65").Value() = "Clemenuless"
$ws.Range("L65").Value() = "Segunda"
$ws.Range("M65").Value() = 300
$ws.Range("N65").Value() = 13000
$ws.Range("O65").Value() = 14000
$ws.Range("P65").Value() = 13500
$ws.Range("Q65").Value() = "$/caja 20 kilos"
$ws.Range("R65").Value() = "Región de Coquimbo"
$ws.Range("S65").Value() = 675
$ws.Range("T65").Value() = 20
